$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tracking Sheet")

# --- Rows 4-6: only the Week (column B) value was blank before; fill it in. ---
$ws.Range("B4").Value = 5
$ws.Range("B5").Value = 5
$ws.Range("B6").Value = 5

# --- New / restyled rows: copy the row-6 cell formatting (dates/body text style) ---
# into each row that needs it, then write the rows values.

# Row 7
$ws.Range("A6:G6").Copy() | Out-Null
$ws.Range("A7:G7").PasteSpecial(-4122) | Out-Null
$ws.Range("A7").Value = 43555
$ws.Range("B7").Value = 5
$ws.Range("C7").Value = 'Android App Development'
$ws.Range("D7").Value = 'Saumil'
$ws.Range("E7").Value = '1. Designed a simple app to discovery bluetooth devices nearby
2. Added buttons to turn ON and OFF bluetooth on the android phone.
3. Displayed the discovered bluetooth devices in a list view'
$ws.Range("F7").Value = 8
$ws.Range("G7").Value = 'BRIDGE'
$ws.Rows.Item(7).RowHeight = 38.25

# Row 8
$ws.Range("A6:G6").Copy() | Out-Null
$ws.Range("A8:G8").PasteSpecial(-4122) | Out-Null
$ws.Range("A8").Value = 43556
$ws.Range("B8").Value = 5
$ws.Range("C8").Value = 'Documentation & Interfacing'
$ws.Range("D8").Value = 'Vidushi '
$ws.Range("E8").Value = '1. Pushed Bridge code to gitlab
2. Hardware Interfaced compass with Sjone Board (5v)'
$ws.Range("F8").Value = 2
$ws.Range("G8").Value = 'BRIDGE & 
GEO'
$ws.Rows.Item(8).RowHeight = 25.5

# Row 9
$ws.Range("A6:G6").Copy() | Out-Null
$ws.Range("A9:G9").PasteSpecial(-4122) | Out-Null
$ws.Range("A9").Value = 43557
$ws.Range("B9").Value = 5
$ws.Range("C9").Value = 'Data Sheet Reviewing '
$ws.Range("D9").Value = 'Vidushi '
$ws.Range("E9").Value = '1. Go through CMPS 11 datasheet in i2c & Serial Mode and read about it register settings and configurations '
$ws.Range("F9").Value = 2
$ws.Range("G9").Value = 'CMPS 11 - GEO'
$ws.Rows.Item(9).RowHeight = 25.5

# Row 10
$ws.Range("A6:G6").Copy() | Out-Null
$ws.Range("A10:G10").PasteSpecial(-4122) | Out-Null
$ws.Range("A10").Value = 43558
$ws.Range("B10").Value = 5
$ws.Range("C10").Value = 'Documentation'
$ws.Range("D10").Value = 'Vidushi '
$ws.Range("E10").Value = '1. Updated wiki report with Project Title details and Car Top & side views.
2. Added Abstract and How self driving car works on Wiki.
3. Added Introduction & onjectives of Tech Savy on Wiki with High Level System Design.'
$ws.Range("F10").Value = 3
$ws.Range("G10").Value = 'WIKI Report'
$ws.Rows.Item(10).RowHeight = 63.75

# Row 11
$ws.Range("A6:G6").Copy() | Out-Null
$ws.Range("A11:G11").PasteSpecial(-4122) | Out-Null
$ws.Range("A11").Value = 43559
$ws.Range("B11").Value = 5
$ws.Range("C11").Value = 'Documentation'
$ws.Range("D11").Value = 'Vidushi '
$ws.Range("E11").Value = '1. Updated wiki report with  Team Members & Technical Responsibilities.
2. Updated wiki report with  Team Members &   Administrative Responsibilities.
3. Updated wiki report with  Team Deliverables Schedule.
4. Added BILL OF MATERIALS (GENERAL PARTS) and Details on the Wiki.
5. Changed the Wiki Font , format and design for our team and added color codes for all modules.'
$ws.Range("F11").Value = 4
$ws.Range("G11").Value = 'WIKI Report'
$ws.Rows.Item(11).RowHeight = 102

# Row 12
$ws.Range("A6:G6").Copy() | Out-Null
$ws.Range("A12:G12").PasteSpecial(-4122) | Out-Null
$ws.Range("A12").Value = 43560
$ws.Range("B12").Value = 5
$ws.Range("C12").Value = 'Interfacing, Coding & Testing'
$ws.Range("D12").Value = 'Vidushi '
$ws.Range("E12").Value = '1. Interfacing of CMPS11 with SJONE board on 3.3V.
2. Interfaced & Implement I2C Mode with CMP11 on SJOne. 
3. Worked on all axis Calibration of CMPS11 using registers 0xF0, 0xF5 and 0xF6.
4. Implemented start & stop caliberation mode for CMPS 11 using command registers on SWITCH.
5. Implemented factory caliberation mode for CMPS11 on Switch in order to revert caliberation.'
$ws.Range("F12").Value = 4
$ws.Range("G12").Value = 'GEO - CMPS11'
$ws.Rows.Item(12).RowHeight = 102

# Row 13
$ws.Range("A6:G6").Copy() | Out-Null
$ws.Range("A13:G13").PasteSpecial(-4122) | Out-Null
$ws.Range("A13").Value = 43561
$ws.Range("B13").Value = 5
$ws.Range("C13").Value = 'Interfacing, Coding & Testing'
$ws.Range("D13").Value = 'Vidushi '
$ws.Range("E13").Value = '1. Study & Implemented I2C Read/Write functions for CMPS11.
2. Implemented Heading Angle Calculation Functionality for CMPS11.
3. Refactoring of code for Geo Controller.
4. Implemented C wrapper for I2C2.cpp.'
$ws.Range("F13").Value = 4
$ws.Range("G13").Value = 'GEO - CMPS11'
$ws.Rows.Item(13).RowHeight = 51

# Row 14
$ws.Range("A6:G6").Copy() | Out-Null
$ws.Range("A14:G14").PasteSpecial(-4122) | Out-Null
$ws.Range("A14").Value = 43563
$ws.Range("B14").Value = 5
$ws.Range("C14").Value = 'Documentation'
$ws.Range("D14").Value = 'Vidushi '
$ws.Range("E14").Value = '1. Updated wiki report with  Team Deliverables Schedule Till Final Demo
2. Updated HW Block Diagrams for Bluetooth and Geo Controller'
$ws.Range("F14").Value = 2
$ws.Range("G14").Value = 'WIKI Report'
$ws.Rows.Item(14).RowHeight = 38.25

# Row 15
$ws.Range("A6:G6").Copy() | Out-Null
$ws.Range("A15:G15").PasteSpecial(-4122) | Out-Null
$ws.Range("A15").Value = 43565
$ws.Range("B15").Value = '4,5,6'
$ws.Range("C15").Value = 'Hardware Design and Motor Control'
$ws.Range("D15").Value = 'Vatsal'
$ws.Range("E15").Value = '1. Tested RC car DC and Servo Motor with the basic PWM API driver and performed ESC calibration before starting the motor
2. Started designing schematics and PCB layout for the RC car project      3. Designed power circuit and identified the components required(BOM)   4. Updated changes to schematic based on the previous project and team''s inputs
5. Helped Jay to develop motor driver implementing CAN and DBC'
$ws.Range("F15").Value = '20 (Combined for week 4,5,6)'
$ws.Range("G15").Value = 'PCB and Motor Controller'
$ws.Rows.Item(15).RowHeight = 89.25

# --- A few cells need the "wrap" style (s=9) instead of the plain body style (s=8); ---
# copy that format from an existing s=9 cell (E6) onto them, then restore their values.
$ws.Range("E6").Copy() | Out-Null
$ws.Range("G8").PasteSpecial(-4122) | Out-Null
$ws.Range("G8").Value = 'BRIDGE & 
GEO'

$ws.Range("E6").Copy() | Out-Null
$ws.Range("F15:G15").PasteSpecial(-4122) | Out-Null
$ws.Range("F15").Value = '20 (Combined for week 4,5,6)'
$ws.Range("G15").Value = 'PCB and Motor Controller'
